$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1) — subject count updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) — updated meanEMG / legmaxROM figures
$ws.Range("B2").Value = 13.754862395879627
$ws.Range("C2").Value = 24.873432544317609
$ws.Range("D2").Value = 31.449211181178271
$ws.Range("E2").Value = 24.797060122536578

# Row 3 (STR) — updated meanEMG / legmaxROM figures
$ws.Range("B3").Value = 10.911190691211516
$ws.Range("C3").Value = 14.383597367489955
$ws.Range("D3").Value = 42.27627787045401
$ws.Range("E3").Value = 17.107705943601673

# Selection now only spans the updated columns (B:E) instead of the full table
$ws.Range("B1:E3").Select() | Out-Null
